# Weekly update: a new Piña (Vega Monumental Concepción) price record was
# reported for the latest week. Insert it as the new row 77 and push the
# existing records (old rows 77-132) down by one (to 78-133).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 77:132 down to 78:133, leaving a blank row 77 to fill in.
$ws.Rows("77").Insert()

# Populate the new row 77 with this week's reported record.
$ws.Range("A77").Value = 11
$ws.Range("B77").Value = "Vega Monumental Concepción"
$ws.Range("C77").Value = "Bíobío"
$ws.Range("D77").Value = 44574
$ws.Range("E77").Value = 8
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100108
$ws.Range("H77").Value = "Tropicales y subtropicales"
$ws.Range("I77").Value = 100108005
$ws.Range("J77").Value = "Piña"
$ws.Range("K77").Value = "Caramelo"
$ws.Range("L77").Value = "Segunda"
$ws.Range("M77").Value = 200
$ws.Range("N77").Value = 16000
$ws.Range("O77").Value = 16500
$ws.Range("P77").Value = 16250
$ws.Range("Q77").Value = "$/caja 14 unidades"
$ws.Range("R77").Value = "Ecuador"
$ws.Range("S77").Value = 1161
$ws.Range("T77").Value = 14
